$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values ("candle" row)
$ws.Range("B2").Value = 1791
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 21492

# Delete row 3 entirely (the "sylvia-breitenberg" product row), shifting rows up
$ws.Rows("3").Delete()

# The grand-total that used to live at E4 is now at E3 (after the delete, old E4 becomes E3)
$ws.Range("E3").Value = 21492

# Update the shared string used for the remaining product name in A2
$ws.Range("A2").Value = "tre-wunsch"

# Column widths (subtract 5/6 to compensate for the engine's internal padding
# so the stored OOXML "width" attribute ends up exactly as intended)
$ws.Columns("A").ColumnWidth = 20 - 5/6
$ws.Columns("B").ColumnWidth = 23 - 5/6
$ws.Columns("C").ColumnWidth = 32 - 5/6
$ws.Columns("D").ColumnWidth = 24 - 5/6
$ws.Columns("E").ColumnWidth = 10 - 5/6
